$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("TargetType"). This shifts the old
# C..J columns (TargetType .. description) one column to the right (D..K),
# carrying their values, shared-string refs and styles with them.
$ws.Columns("C").Insert()

# New column C is "SkillType": a flag for skills that carry an
# attacked/"drawback" behavior. Header + default (0) for the existing
# skills, and 1 for the three new magic skills (rows 23-25).
$ws.Range("C1").Value = "SkillType"
$ws.Range("C2:C22").Value = 0
$ws.Range("C23:C25").Value = 1
$ws.Range("C23:C25").VerticalAlignment = -4108

# Restore the active selection the author left on the sheet.
$ws.Range("C20").Select()

Write-Host "Applied SkillType column insert"
